$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Add the two new rows of triglyceride data to the "exported" sheet (Sheet2).
# Copy the date formatting (style) from an existing date cell so the new
# date cells reuse the same cell style instead of creating a new one.
$ws2.Range("C2").Copy()
$ws2.Range("C8:C9").PasteSpecial(-4122)

$ws2.Range("A8").Value = "Mature Adipocytes"
$ws2.Range("B8").Value = 1450.3205128205129
$ws2.Range("C8").Value = 42972

$ws2.Range("A9").Value = "Mature Adipocytes + Dexamethasone"
$ws2.Range("B9").Value = 590.27777777777783
$ws2.Range("C9").Value = 42972

# Restore each sheet's ruler visibility (turns off the showRuler="0" override
# that was present in the source workbook) and re-establish the selections.
$ws1.Activate()
$ws1.Range("A7:C10").Select()
$excel.ActiveWindow.DisplayRuler = $true

$ws2.Activate()
$ws2.Range("C10").Select()
$excel.ActiveWindow.DisplayRuler = $true
